# Append the newest Adafruit IO reading as row 80 of the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 80

$ws.Cells.Item($newRow, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($newRow, 2).Value = "temperature"

# "25" looks numeric, but the source data stores it as text (like every
# other Value cell in this sheet) - use a leading apostrophe so Excel
# keeps it as text, then restore the default "Normal" style so no stray
# number-format/quote-prefix styling is left on the cell.
$ws.Cells.Item($newRow, 3).Value = "'25"
$ws.Cells.Item($newRow, 3).Style = "Normal"

$ws.Cells.Item($newRow, 4).Value = "N/A"
$ws.Cells.Item($newRow, 5).Value = "N/A"
$ws.Cells.Item($newRow, 6).Value = "N/A"
